$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Metadata sheet (sheet 1): update several property values in place
# ---------------------------------------------------------------------------

# Version: 0.1.6 -> 0.1.7
$ws1.Range("A3").Value = "Version"
$ws1.Range("B3").Value = "0.1.7"

# Status: active -> draft
$ws1.Range("A6").Value = "Status"
$ws1.Range("B6").Value = "draft"

# Date: refreshed timestamp
$ws1.Range("A8").Value = "Date"
$ws1.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Contact (row 10): new publisher contact text (replacing "No display for ContactDetail")
$ws1.Range("A10").Value = "Contact"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11 used to duplicate the old Contact text; it now becomes a second Contact
# entry with a name and e-mail address.
$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# ---------------------------------------------------------------------------
# Insert a new "Jurisdiction" row at position 12, pushing the former rows
# 12-15 (Description, Purpose, Copyright, Immutable) down to rows 13-16.
# ---------------------------------------------------------------------------

# First make sure the new row that will appear at the bottom (row 16) has the
# same formatting (borders/fill/alignment) as the existing data rows.
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift the contents of rows 12-15 down into rows 13-16, working from the
# bottom up so that a row's data is moved before it gets overwritten by the
# row above it. Only values are moved (formatting is already correct for
# every row in this range).
for ($r = 15; $r -ge 12; $r--) {
  $dstRow = $r + 1
  $ws1.Range("A$dstRow`:B$dstRow").ClearContents()
  $ws1.Range("A$r`:B$r").Copy()
  $ws1.Range("A$dstRow`:B$dstRow").PasteSpecial(-4104)
}
$excel.CutCopyMode = $false

# Finally, populate the now-empty row 12 with the new "Jurisdiction" property.
# (Its value is blank, same as the Purpose/Copyright rows.)
$ws1.Cells.Item(12, 1).Value = "Jurisdiction"
$ws1.Cells.Item(12, 2).Value = ""

Write-Output "done"
